# This edit reorders the data rows (2-24) of the sheet: the values in
# columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) are shuffled between
# rows while every other column stays tied to its original row.
#
# Mapping: new row R gets the D/J/K/L/M/P values that used to live in
# old row Map[R] (1-based worksheet row numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row number -> source (old) row number
$map = @{
    2  = 19
    3  = 14
    4  = 5
    5  = 4
    6  = 24
    7  = 2
    8  = 20
    9  = 17
    10 = 12
    11 = 7
    12 = 9
    13 = 16
    14 = 6
    15 = 21
    16 = 8
    17 = 13
    18 = 11
    19 = 10
    20 = 22
    21 = 23
    22 = 15
    23 = 3
    24 = 18
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot the current values of the columns that move, for every row,
# before we start overwriting cells.
$snapshot = @{}
for ($r = 2; $r -le 24; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write back the shuffled values according to the mapping.
for ($r = 2; $r -le 24; $r++) {
    $src = $map[$r]
    $srcVals = $snapshot[$src]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $srcVals[$col]
    }
}
